{"js": "// BitCell.docx edit: add a comma in the Chapter 2 paragraph, and append new\n// story content (a Chapter 3 body paragraph + blank line + \"Chapter 4\"\n// heading) right after the existing \"Chapter 3\" heading paragraph.\n\nconst body = context.document.body;\n\n// --- 1. Insert a comma right after \"...these trees\" and before\n//        \" and the only way to save them was by going to cyberspace.\" ---\nconst tailSearch = body.search(\n  \" and the only way to save them was by going to cyberspace.\",\n  { matchCase: true }\n);\ntailSearch.load(\"items\");\nawait context.sync();\n\nif (tailSearch.items.length > 0) {\n  tailSearch.items[0].insertText(\",\", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// --- 2. Locate the \"Chapter 3\" heading paragraph and insert the new\n//        Chapter 3 story paragraph, a blank paragraph, and a \"Chapter 4\"\n//        heading right after it. ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet chapter3Paragraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"Chapter 3\") {\n    chapter3Paragraph = paragraph;\n    break;\n  }\n}\n\nif (chapter3Paragraph) {\n  const storyText =\n    \"Shean got near one of the trees and put the device close to it. \" +\n    \"In that second, Shean was transported to cyberspace. Cyberspace \" +\n    \"looked like a pixelated house in 3rd person. There was a moving \" +\n    \"object, When Shean saw this creature, it started to shoot at him. \" +\n    \"Quickly he realized that he must shoot them to free the souls. \" +\n    \"After defeating the objects, he was sent back to earth.\";\n\n  const storyParagraph = chapter3Paragraph.insertParagraph(\n    storyText,\n    Word.InsertLocation.after\n  );\n  const blankParagraph = storyParagraph.insertParagraph(\n    \"\",\n    Word.InsertLocation.after\n  );\n  blankParagraph.insertParagraph(\"Chapter 4\", Word.InsertLocation.after);\n\n  await context.sync();\n}\n", "ps1": "# BitCell.docx edit: add a comma in the Chapter 2 paragraph, and append new\n# story content (a Chapter 3 body paragraph + blank line + \"Chapter 4\"\n# heading) right after the existing \"Chapter 3\" heading paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert a comma right after \"...these trees\" and before\n#        \" and the only way to save them was by going to cyberspace.\" ---\n$tail = $d.Content\n$tail.Find.ClearFormatting()\n$tail.Find.Text = \" and the only way to save them was by going to cyberspace.\"\n$found = $tail.Find.Execute()\nif ($found) {\n    $tail.Collapse(1)   # wdCollapseStart\n    $tail.InsertBefore(\",\")\n}\n\n# --- 2. Locate the \"Chapter 3\" heading paragraph and insert the new\n#        Chapter 3 story paragraph, a blank paragraph, and a \"Chapter 4\"\n#        heading right after it. ---\n$storyText = \"Shean got near one of the trees and put the device close to it. In that second, Shean was transported to cyberspace. Cyberspace looked like a pixelated house in 3rd person. There was a moving object, When Shean saw this creature, it started to shoot at him. Quickly he realized that he must shoot them to free the souls. After defeating the objects, he was sent back to earth.\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq \"Chapter 3\") {\n        $p.Range.InsertParagraphAfter()\n\n        $storyPara = $d.Paragraphs($i + 1)\n        $storyPara.Range.InsertAfter($storyText)\n        $storyPara.Range.InsertParagraphAfter()\n\n        $blankPara = $d.Paragraphs($i + 2)\n        $blankPara.Range.InsertParagraphAfter()\n\n        $chapter4Para = $d.Paragraphs($i + 3)\n        $chapter4Para.Range.InsertAfter(\"Chapter 4\")\n\n        break\n    }\n}\n"}
